$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.100.42"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "2.054.09"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'229.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'0.615"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'61.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.76%  "
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").Value = "'14.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").Value = "2.355.46"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "'21.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.16%  "
$ws.Range("D15").Value = "'5.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "'0.757"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").Value = "2.055.69"
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("D18").Value = "38.052.34"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("D20").Value = "'69.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("D22").Value = "'226.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "'165.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "'9.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("D29").Value = "'19.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("D32").Value = "'4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").Value = "'4.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  +8.64%  "
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").Value = "'6.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.41%  "
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "1.518.18"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.13%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'97.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.29%  "
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("E44").Value = "  +2.37%  "
$ws.Range("D45").Value = "'0.0921"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E46").Value = "  +1.29%  "
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").Value = "'2.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").Value = "'7.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").Value = "2.245.17"
$ws.Range("E51").Value = "  +1.69%  "
